# Rename the default "image<N>" labels that Word stamps on the three
# logo pictures living in the document's headers/footers:
#   - header BTec logo:                          image1.jpg -> image2.jpg
#   - footer Pearson logo (first Footers item):   image2.png -> image1.png
#   - footer Pearson logo (second Footers item):  image2.png -> image1.png
#
# The pictures themselves (and their alt text / description) are left
# untouched - only the shape's Name changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($shape, $newName) {
    # Selecting the shape's range first and then reaching the shape via
    # $word.Selection is what reliably lets the rename stick for pictures
    # that live in footers as well as headers.
    [void]$shape.Range.Select()
    $selShape = $word.Selection.InlineShapes.Item(1)
    $selShape.Name = $newName
}

# Header -> BTec logo picture: image1.jpg -> image2.jpg
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            Rename-InlinePicture $shp "image2.jpg"
        }
    }
}

# Footers -> Pearson logo picture (appears twice): image2.png -> image1.png
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -like "*PearsonLogo.png") {
            Rename-InlinePicture $shp "image1.png"
        }
    }
}
